$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title: "SENSOR BASE MOBILE APPLICATION PROJECT" -> "SENSOR BASED MOBILE APPLICATION PROJECT"
#    split into three runs: "SENSOR BASE" / "D" / " MOBILE APPLICATION PROJECT"
# ---------------------------------------------------------------------
$titleFind = $d.Content
$titleFind.Find.Execute("SENSOR BASE MOBILE APPLICATION PROJECT", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$titleStart = $titleFind.Start

$insD = $d.Range($titleStart + 11, $titleStart + 11)
$insD.InsertAfter("D")
$dChar = $d.Range($titleStart + 11, $titleStart + 12)
$dChar.Bold = 1
$dChar.Bold = 0

# ---------------------------------------------------------------------
# 2. Insert "TESTING-" before "UI/UX DOCUMENTATION" and move the
#    "_GoBack" bookmark to sit right after it (Word keeps bookmark
#    names unique, so re-adding it removes the old occurrence).
# ---------------------------------------------------------------------
$uiFind = $d.Content
$uiFind.Find.Execute("UI/UX DOCUMENTATION", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$uiStart = $uiFind.Start

$insTesting = $d.Range($uiStart, $uiStart)
$insTesting.InsertAfter("TESTING-")

$bmRange = $d.Range($uiStart + 8, $uiStart + 8)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------
# 3. Underline the standalone "Problem:" / "Solution:" labels, and
#    split the two full-sentence runs so only the label is underlined.
# ---------------------------------------------------------------------

# First Problem:/Solution: pair - just add underline (already standalone runs)
$scope = $d.Range(0, $d.Content.End)
$scope.Find.Execute("Problem:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$scope.Font.Underline = 1

$scope = $d.Range($scope.End, $d.Content.End)
$scope.Find.Execute("Solution:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$scope.Font.Underline = 1

# Second Problem: pair - the label is merged into the sentence run, split it off
$scope = $d.Range($scope.End, $d.Content.End)
$scope.Find.Execute("Problem: According to some users, it would be cool if the 3D models stay still on the plane where it has appeared", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$problem2Start = $scope.Start
$label = $d.Range($problem2Start, $problem2Start + 8)
$label.Font.Underline = 1

$scope = $d.Range($scope.End, $d.Content.End)
$scope.Find.Execute("Solution: There is no solution for this yet. However, we have done as much as we can to make the 3D models look right on the plane and ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$solution2Start = $scope.Start
$label = $d.Range($solution2Start, $solution2Start + 9)
$label.Font.Underline = 1

# Third Problem:/Solution: pair - just add underline
$scope = $d.Range($scope.End, $d.Content.End)
$scope.Find.Execute("Problem:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$scope.Font.Underline = 1

$scope = $d.Range($scope.End, $d.Content.End)
$scope.Find.Execute("Solution:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$scope.Font.Underline = 1

# ---------------------------------------------------------------------
# 4. Clean up the "How To" proofing split: merge the three runs (with
#    the proofErr gramStart/gramEnd markers around "To") back into a
#    single run by doing an in-place Find & Replace across them.
# ---------------------------------------------------------------------
$scope = $d.Content
$scope.Find.Execute("Under the Play button is the How To Play button which leads the a screen showing instructions to play the game.", $true, $false, $false, $false, $false, $true, 1, $false, "Under the Play button is the How To Play button which leads the a screen showing instructions to play the game.", 2) | Out-Null

$scope = $d.Content
$scope.Find.Execute("Under the How To Play button is the Quit button to exit the application.", $true, $false, $false, $false, $false, $true, 1, $false, "Under the How To Play button is the Quit button to exit the application.", 2) | Out-Null
